$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 34.5
$ws.Range("I11").Value = 34.5
$ws.Range("K11").Value = 34.5
$ws.Range("M11").Value = 105.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1860.3182
$ws.Range("J17").Value = 1860.3182
$ws.Range("L17").Value = 5580.9546
$ws.Range("N17").Value = -5916.9546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 803
$ws.Range("I18").Value = 803
$ws.Range("K18").Value = 803
$ws.Range("M18").Value = -519

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3316.6924
$ws.Range("J70").Value = 4792.4287
$ws.Range("L70").Value = 14377.2861
$ws.Range("N70").Value = -14917.2861

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3316.6924
$ws.Range("J73").Value = 4792.4287
$ws.Range("L73").Value = 14377.2861
$ws.Range("N73").Value = -16249.2861

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 618.25
$ws.Range("I80").Value = 541
$ws.Range("J80").Value = 850
$ws.Range("K80").Value = 1623
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -625
$ws.Range("N80").Value = -4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 618.25
$ws.Range("I83").Value = 541
$ws.Range("J83").Value = 850
$ws.Range("K83").Value = 4869
$ws.Range("L83").Value = 7650
$ws.Range("M83").Value = 123
$ws.Range("N83").Value = -17634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9119.4
$ws.Range("I86").Value = 8532.666999999999
$ws.Range("K86").Value = 8532.666999999999
$ws.Range("M86").Value = -7409.666999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 9119.4
$ws.Range("I89").Value = 8532.666999999999
$ws.Range("K89").Value = 42663.335
$ws.Range("M89").Value = -37047.335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8168.3335
$ws.Range("J113").Value = 4500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -11008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 726.3333
$ws.Range("I127").Value = 726.3333
$ws.Range("K127").Value = 2178.9999
$ws.Range("M127").Value = 2781.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1125.7778
$ws.Range("I129").Value = 891.5
$ws.Range("K129").Value = 2674.5
$ws.Range("M129").Value = 2325.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2547.276
$ws.Range("I137").Value = 1010.8571
$ws.Range("K137").Value = 3032.5713
$ws.Range("M137").Value = -482.5712999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6541.7144
$ws.Range("I74").Value = 6565.6665
$ws.Range("K74").Value = 6565.6665
$ws.Range("M74").Value = -5691.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6541.7144
$ws.Range("I77").Value = 6565.6665
$ws.Range("K77").Value = 32828.3325
$ws.Range("M77").Value = -28460.3325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1579.8
$ws.Range("I110").Value = 1349.75
$ws.Range("K110").Value = 1349.75
$ws.Range("M110").Value = 695.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 294.4
$ws.Range("I94").Value = 291
$ws.Range("J94").Value = 299.5
$ws.Range("K94").Value = 291
$ws.Range("L94").Value = 299.5
$ws.Range("M94").Value = 160
$ws.Range("N94").Value = -1201.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1810.3636
$ws.Range("I99").Value = 1324.1111
$ws.Range("K99").Value = 1324.1111
$ws.Range("M99").Value = 173.8888999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 14418.223
$ws.Range("I103").Value = 10000
$ws.Range("J103").Value = 14970.5
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 14970.5
$ws.Range("M103").Value = -8828
$ws.Range("N103").Value = -17314.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1668.2727
$ws.Range("I105").Value = 1435.4
$ws.Range("J105").Value = 3997
$ws.Range("K105").Value = 1435.4
$ws.Range("L105").Value = 3997
$ws.Range("M105").Value = 311.5999999999999
$ws.Range("N105").Value = -7491

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5201.4443
$ws.Range("I58").Value = 4688.2856
$ws.Range("K58").Value = 4688.2856
$ws.Range("M58").Value = -4485.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5201.4443
$ws.Range("I136").Value = 4688.2856
$ws.Range("K136").Value = 14064.8568
$ws.Range("M136").Value = -11514.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1499.2222
$ws.Range("I5").Value = 1436.625
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 4309.875
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -4197.875
$ws.Range("N5").Value = -6224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 735.2
$ws.Range("I109").Value = 735.2
$ws.Range("K109").Value = 2205.6
$ws.Range("M109").Value = -1165.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 1997.5
$ws.Range("I115").Value = 1997.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 5992.5
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1499.2222
$ws.Range("I135").Value = 1436.625
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 12929.625
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -10394.625
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 39500
$ws.Range("J93").Value = 39500
$ws.Range("L93").Value = 39500
$ws.Range("N93").Value = -43244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2814.5715
$ws.Range("I102").Value = 2617
$ws.Range("K102").Value = 2617
$ws.Range("M102").Value = -995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1603.7646
$ws.Range("I107").Value = 1147.25
$ws.Range("J107").Value = 2699.4
$ws.Range("K107").Value = 1147.25
$ws.Range("L107").Value = 2699.4
$ws.Range("M107").Value = 772.75
$ws.Range("N107").Value = -6539.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2100
$ws.Range("J126").Value = 2100
$ws.Range("L126").Value = 6300
$ws.Range("N126").Value = -11240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 99850
$ws.Range("J134").Value = 99850
$ws.Range("L134").Value = 299550
$ws.Range("N134").Value = -304620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6705.294
$ws.Range("I46").Value = 5666.6665
$ws.Range("J46").Value = 6927.857
$ws.Range("K46").Value = 5666.6665
$ws.Range("L46").Value = 6927.857
$ws.Range("M46").Value = -5478.6665
$ws.Range("N46").Value = -7303.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1215.9166
$ws.Range("J55").Value = 1366.3334
$ws.Range("L55").Value = 1366.3334
$ws.Range("N55").Value = -1712.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7892.857
$ws.Range("J68").Value = 8166.6665
$ws.Range("L68").Value = 8166.6665
$ws.Range("N68").Value = -9664.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 7892.857
$ws.Range("J71").Value = 8166.6665
$ws.Range("L71").Value = 40833.3325
$ws.Range("N71").Value = -48321.3325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2715.739
$ws.Range("I82").Value = 1415.2307
$ws.Range("K82").Value = 1415.2307
$ws.Range("M82").Value = -1054.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2715.739
$ws.Range("I85").Value = 1415.2307
$ws.Range("K85").Value = 1415.2307
$ws.Range("M85").Value = -167.2307000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 798.25
$ws.Range("I93").Value = 798.25
$ws.Range("K93").Value = 798.25
$ws.Range("M93").Value = 449.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8464.571
$ws.Range("J62").Value = 8944.666999999999
$ws.Range("L62").Value = 8944.666999999999
$ws.Range("N62").Value = -10192.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8464.571
$ws.Range("J65").Value = 8944.666999999999
$ws.Range("L65").Value = 44723.335
$ws.Range("N65").Value = -50963.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 931.2
$ws.Range("J100").Value = 1575
$ws.Range("L100").Value = 3150
$ws.Range("N100").Value = -4232
